# Scheduled runner update: refresh currentAveragePrice / Leve price & profit
# columns (H-N) for a batch of leves across multiple job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H112").Value = 1477.2413
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1497.8572
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 4493.571599999999
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -6709.571599999999

$ws.Range("H129").Value = 3437.6924
$ws.Range("I129").Value = 400
$ws.Range("J129").Value = 3690.8333
$ws.Range("K129").Value = 1200
$ws.Range("L129").Value = 11072.4999
$ws.Range("M129").Value = 3800
$ws.Range("N129").Value = -21072.4999

$ws.Range("H132").Value = 1563659.9
$ws.Range("I132").Value = 1954046.8
$ws.Range("J132").Value = 2112.25
$ws.Range("K132").Value = 5862140.4
$ws.Range("L132").Value = 6336.75
$ws.Range("M132").Value = -5859610.4
$ws.Range("N132").Value = -11396.75

$ws = $wb.Worksheets("ARM")
$ws.Range("H33").Value = 1302521.8
$ws.Range("I33").Value = 5000000
$ws.Range("J33").Value = 70029
$ws.Range("K33").Value = 5000000
$ws.Range("L33").Value = 70029
$ws.Range("M33").Value = -4999671
$ws.Range("N33").Value = -70687

$ws.Range("H45").Value = 1137.4231
$ws.Range("I45").Value = 1107.95
$ws.Range("J45").Value = 1235.6666
$ws.Range("K45").Value = 1107.95
$ws.Range("L45").Value = 1235.6666
$ws.Range("M45").Value = -730.95
$ws.Range("N45").Value = -1989.6666

$ws = $wb.Worksheets("BSM")
$ws.Range("H31").Value = 17960
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 22325
$ws.Range("K31").Value = 500
$ws.Range("L31").Value = 22325
$ws.Range("M31").Value = -248
$ws.Range("N31").Value = -22829

$ws.Range("H134").Value = 1705.5
$ws.Range("I134").Value = 1225.9333
$ws.Range("J134").Value = 2733.1428
$ws.Range("K134").Value = 3677.7999
$ws.Range("L134").Value = 8199.428400000001
$ws.Range("M134").Value = -1142.7999
$ws.Range("N134").Value = -13269.4284

$ws = $wb.Worksheets("CRP")
$ws.Range("H3").Value = 170166.67
$ws.Range("I3").Value = 500000
$ws.Range("J3").Value = 5250
$ws.Range("K3").Value = 500000
$ws.Range("L3").Value = 5250
$ws.Range("M3").Value = -499887
$ws.Range("N3").Value = -5476

$ws.Range("H5").Value = 2214.8572
$ws.Range("I5").Value = 149
$ws.Range("J5").Value = 4969.3335
$ws.Range("K5").Value = 149
$ws.Range("L5").Value = 4969.3335
$ws.Range("M5").Value = -37
$ws.Range("N5").Value = -5193.3335

$ws.Range("H11").Value = 12781.2
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 12781.2
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 12781.2
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -13061.2

$ws.Range("H13").Value = 5261.8
$ws.Range("I13").Value = 104
$ws.Range("J13").Value = 6551.25
$ws.Range("K13").Value = 104
$ws.Range("L13").Value = 6551.25
$ws.Range("M13").Value = 35
$ws.Range("N13").Value = -6829.25

$ws.Range("H31").Value = 1880.4651
$ws.Range("I31").Value = 1561.8857
$ws.Range("J31").Value = 3274.25
$ws.Range("K31").Value = 1561.8857
$ws.Range("L31").Value = 3274.25
$ws.Range("M31").Value = -1266.8857
$ws.Range("N31").Value = -3864.25

$ws.Range("H34").Value = 1880.4651
$ws.Range("I34").Value = 1561.8857
$ws.Range("J34").Value = 3274.25
$ws.Range("K34").Value = 1561.8857
$ws.Range("L34").Value = 3274.25
$ws.Range("M34").Value = -1359.8857
$ws.Range("N34").Value = -3678.25

$ws.Range("H36").Value = 44374
$ws.Range("I36").Value = 20048
$ws.Range("J36").Value = 68700
$ws.Range("K36").Value = 20048
$ws.Range("L36").Value = 68700
$ws.Range("M36").Value = -19660
$ws.Range("N36").Value = -69476

$ws.Range("H40").Value = 44374
$ws.Range("I40").Value = 20048
$ws.Range("J40").Value = 68700
$ws.Range("K40").Value = 20048
$ws.Range("L40").Value = 68700
$ws.Range("M40").Value = -19888
$ws.Range("N40").Value = -69020

$ws.Range("H58").Value = 903.9268
$ws.Range("I58").Value = 744.93335
$ws.Range("J58").Value = 1337.5454
$ws.Range("K58").Value = 744.93335
$ws.Range("L58").Value = 1337.5454
$ws.Range("M58").Value = -541.93335
$ws.Range("N58").Value = -1743.5454

$ws.Range("H132").Value = 2755.4443
$ws.Range("I132").Value = 2325
$ws.Range("J132").Value = 3616.3333
$ws.Range("K132").Value = 6975
$ws.Range("L132").Value = 10848.9999
$ws.Range("M132").Value = -4445
$ws.Range("N132").Value = -15908.9999

$ws.Range("H136").Value = 903.9268
$ws.Range("I136").Value = 744.93335
$ws.Range("J136").Value = 1337.5454
$ws.Range("K136").Value = 2234.80005
$ws.Range("L136").Value = 4012.6362
$ws.Range("M136").Value = 315.1999500000002
$ws.Range("N136").Value = -9112.636200000001

$ws = $wb.Worksheets("CUL")
$ws.Range("H92").Value = 1585.1428
$ws.Range("I92").Value = 1374
$ws.Range("J92").Value = 1866.6666
$ws.Range("K92").Value = 4122
$ws.Range("L92").Value = 5599.9998
$ws.Range("M92").Value = -2874
$ws.Range("N92").Value = -8095.9998

$ws.Range("H98").Value = 167.33333
$ws.Range("I98").Value = 167.33333
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 501.99999
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 996.00001
$ws.Range("N98").ClearContents()

$ws.Range("H104").Value = 1364.5264
$ws.Range("I104").Value = 855
$ws.Range("J104").Value = 1424.4706
$ws.Range("K104").Value = 2565
$ws.Range("L104").Value = 4273.4118
$ws.Range("M104").Value = 56
$ws.Range("N104").Value = -9515.4118

$ws.Range("H105").Value = 8000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 8000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 24000
$ws.Range("N105").Value = -29242

$ws = $wb.Worksheets("GSM")
$ws.Range("H97").Value = 657.8570999999999
$ws.Range("I97").Value = 488.66666
$ws.Range("J97").Value = 962.4
$ws.Range("K97").Value = 488.66666
$ws.Range("L97").Value = 962.4
$ws.Range("M97").Value = 7.333340000000021
$ws.Range("N97").Value = -1954.4

$ws = $wb.Worksheets("LTW")
$ws.Range("H82").Value = 7501.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 7501.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 7501.5
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -8223.5

$ws.Range("H85").Value = 7501.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 7501.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 7501.5
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -9997.5
